$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 4002
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 4002
$ws.Range("M19").Value = ""
$ws.Range("N19").Value = -4352
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 4002
$ws.Range("H40").Value = 3442.5715
$ws.Range("H69").Value = 28561.625
$ws.Range("I69").Value = 17699.6
$ws.Range("K69").Value = 53098.8
$ws.Range("M69").Value = -52224.8
$ws.Range("H72").Value = 28561.625
$ws.Range("I72").Value = 17699.6
$ws.Range("K72").Value = 159296.4
$ws.Range("M72").Value = -154928.4
$ws.Range("H133").Value = 86399.60000000001
$ws.Range("J133").Value = 86399.60000000001
$ws.Range("N133").Value = -96519.60000000001
$ws.Range("L133").Value = 86399.60000000001
$ws.Range("H135").Value = 1478.2273
$ws.Range("I135").Value = 1514.95
$ws.Range("K135").Value = 13634.55
$ws.Range("M135").Value = -11099.55
$ws.Range("H137").Value = 3839.238
$ws.Range("J137").Value = 5859.5
$ws.Range("L137").Value = 17578.5
$ws.Range("N137").Value = -22678.5
$ws.Range("L138").Value = 9568.0905
$ws.Range("H138").Value = 2746.2385
$ws.Range("I138").Value = 1416.8636
$ws.Range("J138").Value = 3189.3635
$ws.Range("M138").Value = 889.4092000000001
$ws.Range("N138").Value = -19848.0905
$ws.Range("K138").Value = 4250.5908
$ws.Range("H141").Value = 2344.7896
$ws.Range("I141").Value = 2249.2942
$ws.Range("K141").Value = 6747.882599999999
$ws.Range("M141").Value = -1567.882599999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 14323002
$ws.Range("I61").Value = 25005238
$ws.Range("K61").Value = 25005238
$ws.Range("M61").Value = -25005026
$ws.Range("H74").Value = 9268507
$ws.Range("I74").Value = 14707941
$ws.Range("K74").Value = 14707941
$ws.Range("M74").Value = -14707067
$ws.Range("H77").Value = 9268507
$ws.Range("I77").Value = 14707941
$ws.Range("K77").Value = 73539705
$ws.Range("M77").Value = -73535337
$ws.Range("H94").Value = 47963
$ws.Range("J94").Value = 47963
$ws.Range("L94").Value = 47963
$ws.Range("N94").Value = -49765
$ws.Range("M97").Value = -828.625
$ws.Range("L97").Value = 2000
$ws.Range("H97").Value = 1364.3529
$ws.Range("I97").Value = 1324.625
$ws.Range("J97").Value = 2000
$ws.Range("K97").Value = 1324.625
$ws.Range("N97").Value = -2992
$ws.Range("H106").Value = 49980
$ws.Range("J106").Value = 49980
$ws.Range("L106").Value = 49980
$ws.Range("N106").Value = -52504
$ws.Range("H132").Value = 7466.7856
$ws.Range("I132").Value = 2539.75
$ws.Range("J132").Value = 14036.167
$ws.Range("K132").Value = 7619.25
$ws.Range("M132").Value = -5089.25
$ws.Range("N132").Value = -47168.501
$ws.Range("L132").Value = 42108.501
$ws.Range("H136").Value = 14323002
$ws.Range("I136").Value = 25005238
$ws.Range("K136").Value = 75015714
$ws.Range("M136").Value = -75013164

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("I20").Value = 6269.2
$ws.Range("H20").Value = 6269.2
$ws.Range("K20").Value = 6269.2
$ws.Range("M20").Value = -6022.2
$ws.Range("I134").Value = 1845.8889
$ws.Range("H134").Value = 44153.543
$ws.Range("K134").Value = 5537.6667
$ws.Range("M134").Value = -3002.6667
$ws.Range("L139").Value = 75000
$ws.Range("H139").Value = 68750
$ws.Range("I139").Value = 50000
$ws.Range("J139").Value = 75000
$ws.Range("M139").Value = -44860
$ws.Range("N139").Value = -85280
$ws.Range("K139").Value = 50000

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I31").Value = 22889.75
$ws.Range("H31").Value = 909014.1
$ws.Range("K31").Value = 22889.75
$ws.Range("M31").Value = -22594.75
$ws.Range("H34").Value = 909014.1
$ws.Range("I34").Value = 22889.75
$ws.Range("K34").Value = 22889.75
$ws.Range("M34").Value = -22687.75
$ws.Range("H51").Value = 42600
$ws.Range("H53").Value = 39777.168
$ws.Range("J53").Value = 39777.168
$ws.Range("N53").Value = -40991.168
$ws.Range("L53").Value = 39777.168
$ws.Range("H61").Value = 42600
$ws.Range("L99").Value = 3073.75
$ws.Range("H99").Value = 3330.875
$ws.Range("I99").Value = 3416.5833
$ws.Range("J99").Value = 3073.75
$ws.Range("M99").Value = -1918.5833
$ws.Range("N99").Value = -6069.75
$ws.Range("K99").Value = 3416.5833
$ws.Range("H107").Value = 725.6667
$ws.Range("I107").Value = 571.375
$ws.Range("K107").Value = 571.375
$ws.Range("M107").Value = 1348.625
$ws.Range("L126").Value = 9221.25
$ws.Range("H126").Value = 3330.875
$ws.Range("I126").Value = 3416.5833
$ws.Range("J126").Value = 3073.75
$ws.Range("M126").Value = -7779.749899999999
$ws.Range("N126").Value = -14161.25
$ws.Range("K126").Value = 10249.7499
$ws.Range("H132").Value = 1744.4464
$ws.Range("I132").Value = 1422.102
$ws.Range("J132").Value = 4000.8572
$ws.Range("K132").Value = 4266.306
$ws.Range("M132").Value = -1736.306
$ws.Range("N132").Value = -17062.5716
$ws.Range("L132").Value = 12002.5716

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1627.4166
$ws.Range("J113").Value = 1792.8
$ws.Range("L113").Value = 5378.4
$ws.Range("N113").Value = -9718.4
$ws.Range("H132").Value = 2382.818
$ws.Range("J132").Value = 2470.3333
$ws.Range("L132").Value = 22232.9997
$ws.Range("N132").Value = -27292.9997

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2916.375
$ws.Range("I97").Value = 2904.4285
$ws.Range("K97").Value = 2904.4285
$ws.Range("M97").Value = -2408.4285
$ws.Range("H126").Value = 4797.5
$ws.Range("J126").Value = 4795
$ws.Range("L126").Value = 14385
$ws.Range("N126").Value = -19325
$ws.Range("H132").Value = 166705010
$ws.Range("I132").Value = 333340000
$ws.Range("J132").Value = 70010.664
$ws.Range("K132").Value = 1000020000
$ws.Range("M132").Value = -1000017470
$ws.Range("N132").Value = -215091.992
$ws.Range("L132").Value = 210031.992

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 62333.35
$ws.Range("I7").Value = 2174.875
$ws.Range("J7").Value = 115807.555
$ws.Range("K7").Value = 2174.875
$ws.Range("M7").Value = -2062.875
$ws.Range("N7").Value = -116031.555
$ws.Range("L7").Value = 115807.555
$ws.Range("H40").Value = 4601.353
$ws.Range("I40").Value = 4026.5
$ws.Range("K40").Value = 4026.5
$ws.Range("M40").Value = -3890.5
$ws.Range("H82").Value = 1139.1538
$ws.Range("I82").Value = 604.375
$ws.Range("K82").Value = 604.375
$ws.Range("M82").Value = -243.375
$ws.Range("H85").Value = 1139.1538
$ws.Range("I85").Value = 604.375
$ws.Range("K85").Value = 604.375
$ws.Range("M85").Value = 643.625
$ws.Range("H93").Value = 71437900
$ws.Range("I93").Value = 76932980
$ws.Range("K93").Value = 76932980
$ws.Range("M93").Value = -76931732
$ws.Range("H103").Value = 51926
$ws.Range("J103").Value = 51926
$ws.Range("N103").Value = -54270
$ws.Range("L103").Value = 51926
$ws.Range("H122").Value = 7370.4165
$ws.Range("I122").Value = 7271.778
$ws.Range("J122").Value = 7666.3335
$ws.Range("M122").Value = -19365.334
$ws.Range("K122").Value = 21815.334
$ws.Range("N122").Value = -27899.0005
$ws.Range("L122").Value = 22999.0005
$ws.Range("L126").Value = 347422.665
$ws.Range("H126").Value = 62333.35
$ws.Range("I126").Value = 2174.875
$ws.Range("J126").Value = 115807.555
$ws.Range("M126").Value = -4054.625
$ws.Range("N126").Value = -352362.665
$ws.Range("K126").Value = 6524.625
$ws.Range("H132").Value = 35240.15
$ws.Range("I132").Value = 7203.087
$ws.Range("K132").Value = 21609.261
$ws.Range("M132").Value = -19079.261
$ws.Range("H136").Value = 52867.44
$ws.Range("I136").Value = 7417.6875
$ws.Range("K136").Value = 22253.0625
$ws.Range("M136").Value = -19703.0625

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 24625
$ws.Range("J74").Value = 24625
$ws.Range("N74").Value = -26497
$ws.Range("L74").Value = 24625
$ws.Range("H77").Value = 24625
$ws.Range("J77").Value = 24625
$ws.Range("N77").Value = -83235
$ws.Range("L77").Value = 73875
$ws.Range("H120").Value = 61000
$ws.Range("J120").Value = 61000
$ws.Range("L120").Value = 61000
$ws.Range("N120").Value = -70676
$ws.Range("H122").Value = 16423.812
$ws.Range("I122").Value = 13565.417
$ws.Range("J122").Value = 24999
$ws.Range("M122").Value = -38246.251
$ws.Range("K122").Value = 40696.251
$ws.Range("N122").Value = -79897
$ws.Range("L122").Value = 74997
$ws.Range("H126").Value = 4463.857
$ws.Range("I126").Value = 4169.4
$ws.Range("K126").Value = 12508.2
$ws.Range("M126").Value = -10038.2
$ws.Range("H132").Value = 2848.3215
$ws.Range("I132").Value = 1947.3914
$ws.Range("J132").Value = 6992.6
$ws.Range("K132").Value = 5842.174199999999
$ws.Range("M132").Value = -3312.174199999999
$ws.Range("N132").Value = -26037.8
$ws.Range("L132").Value = 20977.8
$ws.Range("L136").Value = 1042.5
$ws.Range("H136").Value = 790.6316
$ws.Range("I136").Value = 842.7646999999999
$ws.Range("J136").Value = 347.5
$ws.Range("M136").Value = 21.70589999999993
$ws.Range("K136").Value = 2528.2941
$ws.Range("N136").Value = -6142.5
